$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''25.938.78'
$ws.Range("E2").Value = '  -0.54%  '

$ws.Range("D3").Value = '''1.618.72'
$ws.Range("E3").Value = '  -1.27%  '

$ws.Range("E4").Value = '  -0.64%  '

$ws.Range("D5").Value = '''212.84'
$ws.Range("E5").Value = '  -0.79%  '

$ws.Range("D6").Value = '''0.500'
$ws.Range("E6").Value = '  -0.97%  '

$ws.Range("E7").Value = '  -0.69%  '

$ws.Range("E8").Value = '  -0.35%  '

$ws.Range("D9").Value = '''0.0617'
$ws.Range("E9").Value = '  -0.99%  '

$ws.Range("D10").Value = '''18.36'
$ws.Range("E10").Value = '  -1.09%  '

$ws.Range("D11").Value = '''0.0791'
$ws.Range("E11").Value = '  -0.45%  '

$ws.Range("D12").Value = '''1.843.42'
$ws.Range("E12").Value = '  -1.27%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''4.13'
$ws.Range("E13").Value = '  -1.61%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '''1.604.54'
$ws.Range("E14").Value = '  -5.14%  '

$ws.Range("D15").Value = '''0.523'
$ws.Range("E15").Value = '  -1.15%  '

$ws.Range("D16").Value = '''25.925.35'
$ws.Range("E16").Value = '  -0.61%  '

$ws.Range("D17").Value = '''61.54'
$ws.Range("E17").Value = '  -1.27%  '

$ws.Range("D18").Value = '''0.0₃0737'
$ws.Range("E18").Value = '  -1.56%  '

$ws.Range("E19").Value = '  -0.61%  '

$ws.Range("D20").Value = '''191.75'
$ws.Range("E20").Value = '  +0.80%  '

$ws.Range("D21").Value = '''4.25'
$ws.Range("E21").Value = '  -0.54%  '

$ws.Range("D22").Value = '''9.51'
$ws.Range("E22").Value = '  -0.48%  '

$ws.Range("D23").Value = '''6.03'
$ws.Range("E23").Value = '  -1.60%  '

$ws.Range("E24").Value = '  +2.78%  '

$ws.Range("D25").Value = '''143.79'
$ws.Range("E25").Value = '  -0.29%  '

$ws.Range("E26").Value = '  -0.61%  '

$ws.Range("D27").Value = '''1.71'
$ws.Range("E27").Value = '  -2.96%  '

$ws.Range("D28").Value = '''6.64'

$ws.Range("D29").Value = '''15.25'
$ws.Range("E29").Value = '  +0.23%  '

$ws.Range("E30").Value = '  -0.92%  '

$ws.Range("D31").Value = '''0.0478'
$ws.Range("E31").Value = '  -1.32%  '

$ws.Range("E32").Value = '  -1.56%  '

$ws.Range("D33").Value = '''3.10'
$ws.Range("E33").Value = '  -2.49%  '

$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").Value = '''1.50'
$ws.Range("E34").Value = '  -0.62%  '

$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '''2.41'
$ws.Range("E35").Value = '  -1.49%  '

$ws.Range("D36").Value = '''1.126.39'
$ws.Range("E36").Value = '  +0.28%  '

$ws.Range("D37").Value = '''0.842'
$ws.Range("E37").Value = '  -4.32%  '

$ws.Range("E38").Value = '  -3.50%  '

$ws.Range("D39").Value = '''0.512'
$ws.Range("E39").Value = '  -2.02%  '

$ws.Range("E40").Value = '  -1.03%  '

$ws.Range("D41").Value = '''97.86'
$ws.Range("E41").Value = '  -0.87%  '

$ws.Range("D42").Value = '''1.753.01'
$ws.Range("E42").Value = '  -1.17%  '

$ws.Range("E43").Value = '  -4.82%  '

$ws.Range("D44").Value = '''5.07'
$ws.Range("E44").Value = '  -4.13%  '

$ws.Range("D45").Value = '''0.0₆0112'
$ws.Range("E45").Value = '  -2.10%  '

$ws.Range("E46").Value = '  +3.10%  '

$ws.Range("D47").Value = '''54.11'
$ws.Range("E47").Value = '  -1.98%  '

$ws.Range("D48").Value = '''0.0517'
$ws.Range("E48").Value = '  -0.75%  '

$ws.Range("D49").Value = '''0.411'
$ws.Range("E49").Value = '  -1.11%  '

$ws.Range("D50").Value = '''7.49'
$ws.Range("E50").Value = '  -1.01%  '

$ws.Range("E51").Value = '  -0.59%  '
